$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" column (G) for both data rows.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 08:32:27"
$wsOverview.Range("G3").Value = "2016-08-31 08:32:27"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# "Priority" column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# "Correspond Handoff Datetime" column (H)
$wsZhCn.Range("H2").Value = "2016-08-31 08:32:14"
$wsZhCn.Range("H3").Value = "2016-08-31 08:32:14"
# "Correspond Handback DateTime" column (K)
$wsZhCn.Range("K2").Value = "2016-08-31 08:33:19"
$wsZhCn.Range("K3").Value = "2016-08-31 08:33:19"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# "Priority" column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# "Correspond Handoff Datetime" column (H) - shares the Overview's date string
$wsDeDe.Range("H2").Value = "2016-08-31 08:32:27"
$wsDeDe.Range("H3").Value = "2016-08-31 08:32:27"
# "Correspond Handback DateTime" column (K)
$wsDeDe.Range("K2").Value = "2016-08-31 08:33:38"
$wsDeDe.Range("K3").Value = "2016-08-31 08:33:38"
